# Updated cryptos list on Mon Feb 26 13:49:23 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# re-applies the ranking re-sort that swapped ShibaInu/InternetComputer(DFINITY)
# (rows 21-22) and Mantle/WOONetwork (rows 50-51).
#
# Note: several prices (e.g. "393.47") look like plain numbers to Excel's
# auto-detection, but the source data stores them as literal text. We force
# text entry (NumberFormat "@" while writing, then restore the Normal style)
# so the cell keeps the exact literal string instead of being silently
# converted to a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.319.55"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.066.79"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "393.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "3.556.00"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "3.080.23"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "51.336.21"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0487"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.289"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").Value = "2.067.57"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.896"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.71%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.515"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.43%  "
